$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(1,9,5,6,4,2,7,3,8),
    @(6,2,7,1,8,3,9,4,5),
    @(3,4,8,5,7,9,2,1,6),
    @(7,8,9,4,6,5,1,2,3),
    @(4,6,2,3,1,7,8,5,9),
    @(5,3,1,2,9,8,6,7,4),
    @(8,1,4,7,5,6,3,9,2),
    @(9,7,3,8,2,4,5,6,1),
    @(2,5,6,9,3,1,4,8,7)
)

for ($r = 0; $r -lt 9; $r++) {
    for ($c = 0; $c -lt 9; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}
